{"js": "// Update the date paragraph (first paragraph of the document)\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2023-08-22 Tuesday\", \"Replace\");\n\n// Update the division-problem table cells (addressed by row/column so duplicate\n// cell values elsewhere in the table can never cause a wrong cell to be edited)\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  { row: 0, col: 0, text: \"68\u00f77=9, 5\" }, // was \"39\u00f74=9, 3\"\n  { row: 0, col: 1, text: \"15\u00f74=3, 3\" }, // was \"92\u00f77=13, 1\"\n  { row: 0, col: 2, text: \"49\u00f77=7, 0\" }, // was \"45\u00f76=7, 3\"\n  { row: 0, col: 3, text: \"69\u00f79=7, 6\" }, // was \"42\u00f73=14, 0\"\n  { row: 0, col: 4, text: \"63\u00f76=10, 3\" }, // was \"86\u00f73=28, 2\"\n  { row: 4, col: 0, text: \"20\u00f78=2, 4\" }, // was \"26\u00f74=6, 2\"\n  { row: 4, col: 1, text: \"71\u00f78=8, 7\" }, // was \"19\u00f78=2, 3\"\n  { row: 4, col: 2, text: \"24\u00f78=3, 0\" }, // was \"21\u00f79=2, 3\"\n  { row: 4, col: 3, text: \"73\u00f79=8, 1\" }, // was \"39\u00f75=7, 4\"\n  { row: 4, col: 4, text: \"23\u00f73=7, 2\" }, // was \"40\u00f76=6, 4\"\n  { row: 8, col: 0, text: \"96\u00f76=16, 0\" }, // was \"81\u00f76=13, 3\"\n  { row: 8, col: 1, text: \"28\u00f78=3, 4\" }, // was \"73\u00f77=10, 3\"\n  { row: 8, col: 2, text: \"49\u00f76=8, 1\" }, // was \"30\u00f79=3, 3\"\n  { row: 8, col: 3, text: \"92\u00f75=18, 2\" }, // was \"93\u00f76=15, 3\"\n  { row: 8, col: 4, text: \"67\u00f76=11, 1\" }, // was \"75\u00f79=8, 3\"\n  { row: 12, col: 0, text: \"69\u00f74=17, 1\" }, // was \"49\u00f76=8, 1\"\n  { row: 12, col: 1, text: \"41\u00f72=20, 1\" }, // was \"83\u00f72=41, 1\"\n  { row: 12, col: 2, text: \"70\u00f78=8, 6\" }, // was \"95\u00f74=23, 3\"\n  { row: 12, col: 3, text: \"45\u00f73=15, 0\" }, // was \"78\u00f72=39, 0\"\n  { row: 12, col: 4, text: \"13\u00f79=1, 4\" }, // was \"74\u00f73=24, 2\"\n  { row: 16, col: 0, text: \"91\u00f78=11, 3\" }, // was \"56\u00f79=6, 2\"\n  { row: 16, col: 1, text: \"40\u00f78=5, 0\" }, // was \"85\u00f77=12, 1\"\n  { row: 16, col: 2, text: \"83\u00f73=27, 2\" }, // was \"57\u00f76=9, 3\"\n  { row: 16, col: 3, text: \"93\u00f74=23, 1\" }, // was \"16\u00f77=2, 2\"\n  { row: 16, col: 4, text: \"60\u00f75=12, 0\" }, // was \"99\u00f75=19, 4\"\n];\n\nfor (const { row, col, text } of cellUpdates) {\n  const cell = table.getCell(row, col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraph.getRange().insertText(text, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date paragraph (first paragraph of the document)\n$d.Paragraphs.First.Range.Text = \"2023-08-22 Tuesday\"\n\n# Update the division-problem table cells. Cells are addressed by their\n# fixed (row, column) position rather than by searching for their old text,\n# because several of the new values happen to equal *other* cells' old\n# values (e.g. row 9 col 3 becomes \"49\u00f76=8, 1\", which is row 13 col 1's\n# original value) -- a blind Find/Replace could hit the wrong occurrence\n# once earlier replacements have run. Setting Range.Text on a specific\n# Cell is unambiguous and preserves the existing run/paragraph formatting.\n$tbl = $d.Tables.Item(1)\n$tbl.Cell(1, 1).Range.Text = \"68\u00f77=9, 5\"   # was \"39\u00f74=9, 3\"\n$tbl.Cell(1, 2).Range.Text = \"15\u00f74=3, 3\"   # was \"92\u00f77=13, 1\"\n$tbl.Cell(1, 3).Range.Text = \"49\u00f77=7, 0\"   # was \"45\u00f76=7, 3\"\n$tbl.Cell(1, 4).Range.Text = \"69\u00f79=7, 6\"   # was \"42\u00f73=14, 0\"\n$tbl.Cell(1, 5).Range.Text = \"63\u00f76=10, 3\"  # was \"86\u00f73=28, 2\"\n$tbl.Cell(5, 1).Range.Text = \"20\u00f78=2, 4\"   # was \"26\u00f74=6, 2\"\n$tbl.Cell(5, 2).Range.Text = \"71\u00f78=8, 7\"   # was \"19\u00f78=2, 3\"\n$tbl.Cell(5, 3).Range.Text = \"24\u00f78=3, 0\"   # was \"21\u00f79=2, 3\"\n$tbl.Cell(5, 4).Range.Text = \"73\u00f79=8, 1\"   # was \"39\u00f75=7, 4\"\n$tbl.Cell(5, 5).Range.Text = \"23\u00f73=7, 2\"   # was \"40\u00f76=6, 4\"\n$tbl.Cell(9, 1).Range.Text = \"96\u00f76=16, 0\"  # was \"81\u00f76=13, 3\"\n$tbl.Cell(9, 2).Range.Text = \"28\u00f78=3, 4\"   # was \"73\u00f77=10, 3\"\n$tbl.Cell(9, 3).Range.Text = \"49\u00f76=8, 1\"   # was \"30\u00f79=3, 3\"\n$tbl.Cell(9, 4).Range.Text = \"92\u00f75=18, 2\"  # was \"93\u00f76=15, 3\"\n$tbl.Cell(9, 5).Range.Text = \"67\u00f76=11, 1\"  # was \"75\u00f79=8, 3\"\n$tbl.Cell(13, 1).Range.Text = \"69\u00f74=17, 1\" # was \"49\u00f76=8, 1\"\n$tbl.Cell(13, 2).Range.Text = \"41\u00f72=20, 1\" # was \"83\u00f72=41, 1\"\n$tbl.Cell(13, 3).Range.Text = \"70\u00f78=8, 6\"  # was \"95\u00f74=23, 3\"\n$tbl.Cell(13, 4).Range.Text = \"45\u00f73=15, 0\" # was \"78\u00f72=39, 0\"\n$tbl.Cell(13, 5).Range.Text = \"13\u00f79=1, 4\"  # was \"74\u00f73=24, 2\"\n$tbl.Cell(17, 1).Range.Text = \"91\u00f78=11, 3\" # was \"56\u00f79=6, 2\"\n$tbl.Cell(17, 2).Range.Text = \"40\u00f78=5, 0\"  # was \"85\u00f77=12, 1\"\n$tbl.Cell(17, 3).Range.Text = \"83\u00f73=27, 2\" # was \"57\u00f76=9, 3\"\n$tbl.Cell(17, 4).Range.Text = \"93\u00f74=23, 1\" # was \"16\u00f77=2, 2\"\n$tbl.Cell(17, 5).Range.Text = \"60\u00f75=12, 0\" # was \"99\u00f75=19, 4\"\n"}
